$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "20160406_104313"
$ws.Range("B12").Value = 1142.736
$ws.Range("C12").Value = "trim `"space`" and `",`", remove multiple spaces, convert to lower, convert unicode to ascii"
$ws.Range("D12").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E12").Value = "Neuron Network"
$ws.Range("F12").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G12").Value = 0.999333333333333
$ws.Range("H12").Value = 0.920792079207921
$ws.Range("I12").Value = "0 filters: "
$ws.Range("J12").Value = 0.04

$ws.Range("A13").Value = "20160406_110216"
$ws.Range("B13").Value = 1164.935
$ws.Range("C13").Value = "trim `"space`" and `",`", remove multiple spaces, convert to lower, convert unicode to ascii"
$ws.Range("D13").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E13").Value = "Neuron Network"
$ws.Range("F13").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 0.894389438943894
$ws.Range("I13").Value = "0 filters: "
$ws.Range("J13").Value = 0.0447761194029851

$ws.Range("A14").Value = "20160406_112141"
$ws.Range("B14").Value = 1218.655
$ws.Range("C14").Value = "trim `"space`" and `",`", remove multiple spaces, convert to lower, convert unicode to ascii"
$ws.Range("D14").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E14").Value = "Neuron Network"
$ws.Range("F14").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G14").Value = 0.998666666666667
$ws.Range("H14").Value = 0.937293729372937
$ws.Range("I14").Value = "0 filters: "
$ws.Range("J14").Value = 0.0375

$ws.Range("A15").Value = "20160406_114200"
$ws.Range("B15").Value = 1232.143
$ws.Range("C15").Value = "trim `"space`" and `",`", remove multiple spaces, convert to lower, convert unicode to ascii"
$ws.Range("D15").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E15").Value = "Neuron Network"
$ws.Range("F15").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0.884488448844885
$ws.Range("I15").Value = "0 filters: "
$ws.Range("J15").Value = 0.0625

$ws.Range("A16").Value = "20160406_120232"
$ws.Range("B16").Value = 1322.821
$ws.Range("C16").Value = "trim `"space`" and `",`", remove multiple spaces, convert to lower, convert unicode to ascii"
$ws.Range("D16").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E16").Value = "Neuron Network"
$ws.Range("F16").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 0.874587458745875
$ws.Range("I16").Value = "0 filters: "
$ws.Range("J16").Value = 0.0491803278688525

$ws.Range("A17").Value = "20160406_133424"
$ws.Range("B17").Value = 2583.971
$ws.Range("C17").Value = "remove multiple spaces, convert unicode to ascii, trim `"space`" and `",`", convert to lower"
$ws.Range("D17").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E17").Value = "Neuron Network"
$ws.Range("F17").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G17").Value = 0.999333333333333
$ws.Range("H17").Value = 0.900990099009901
$ws.Range("I17").Value = "0 filters: "
$ws.Range("J17").Value = 0.0579710144927536

$ws.Range("A18").Value = "20160406_141728"
$ws.Range("B18").Value = 2664.168
$ws.Range("C18").Value = "remove multiple spaces, convert unicode to ascii, trim `"space`" and `",`", convert to lower"
$ws.Range("D18").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E18").Value = "Neuron Network"
$ws.Range("F18").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 0.897689768976898
$ws.Range("I18").Value = "0 filters: "
$ws.Range("J18").Value = 0.0588235294117647

$ws.Range("A19").Value = "20160406_150152"
$ws.Range("B19").Value = 1758.119
$ws.Range("C19").Value = "remove multiple spaces, convert unicode to ascii, trim `"space`" and `",`", convert to lower"
$ws.Range("D19").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E19").Value = "Neuron Network"
$ws.Range("F19").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G19").Value = 0.999333333333333
$ws.Range("H19").Value = 0.910891089108911
$ws.Range("I19").Value = "0 filters: "
$ws.Range("J19").Value = 0.0416666666666667

$ws.Range("A20").Value = "20160406_153110"
$ws.Range("B20").Value = 1661.828
$ws.Range("C20").Value = "remove multiple spaces, convert unicode to ascii, trim `"space`" and `",`", convert to lower"
$ws.Range("D20").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E20").Value = "Neuron Network"
$ws.Range("F20").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G20").Value = 0.999333333333333
$ws.Range("H20").Value = 0.900990099009901
$ws.Range("I20").Value = "0 filters: "
$ws.Range("J20").Value = 0.0434782608695652

$ws.Range("A21").Value = "20160406_155852"
$ws.Range("B21").Value = 1516.028
$ws.Range("C21").Value = "remove multiple spaces, convert unicode to ascii, trim `"space`" and `",`", convert to lower"
$ws.Range("D21").Value = "8 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2, first_character_type, #`"space`""
$ws.Range("E21").Value = "Neuron Network"
$ws.Range("F21").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 0.874587458745875
$ws.Range("I21").Value = "0 filters: "
$ws.Range("J21").Value = 0.0819672131147541

Write-Host "Rows 12-21 written"
